# The "reviews_count" column (column E) is removed from the sheet.
# Deleting the entire column shifts reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link and latest_review_date one column to
# the left (F->E, G->F, H->G, I->H, J->I, K->J), matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E:E").EntireColumn.Delete()
